$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "Kvareli"
$ws.Name = "Kvareli"

# Row 6 ("Urban"): columns B through M become "..." (N and O already show the
# existing ellipsis character and are left untouched)
$ws.Range("B6:M6").Value = "..."

# Row 7 ("Rural"): columns B through L become "..."; M7 keeps its numeric
# value (4); N7/O7 already show the existing ellipsis character and are left
# untouched
$ws.Range("B7:L7").Value = "..."

# Remove the empty row 8 gap so the old row 9 (the footnote row) becomes row 8
$ws.Rows.Item(8).Delete()
